$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, pushing existing rows 71-88 down to 72-89.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new weekly data record.
$ws.Range("A71").Value = 5
$ws.Range("B71").Value = "Macroferia Regional de Talca"
$ws.Range("C71").Value = "Maule"
$ws.Range("D71").Value = 44900
$ws.Range("E71").Value = 7
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100101
$ws.Range("H71").Value = "Berries"
$ws.Range("I71").Value = 100101001
$ws.Range("J71").Value = "Arándano (blue)"
$ws.Range("K71").Value = "Sin especificar"
$ws.Range("L71").Value = "Primera"
$ws.Range("M71").Value = 270
$ws.Range("N71").Value = 3000
$ws.Range("O71").Value = 3000
$ws.Range("P71").Value = 3000
$ws.Range("Q71").Value = "$/bandeja 2 kilos"
$ws.Range("R71").Value = "Provincia de Curicó"
$ws.Range("S71").Value = 1500
$ws.Range("T71").Value = 2
